$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5078
$ws.Range("C2").Value = 0.64108
$ws.Range("D2").Value = 0.22344
$ws.Range("E2").Value = 0.19839
$ws.Range("F2").Value = "div(div(div(sqrt(sqrt(sqrt(mul(PHS, mul(log(sub(chi, C_m)), r))))), exp(C_0)), exp(C_0)), exp(C_0))"

$wb.Save()
